# Update countries & provincias Spain
# Applies the data refresh described in the commit:
# - Shared-string country order changes (3 countries moved up in the
#   sorted-by-cases list: Somalia, Paraguay, Curazao), which shifts the
#   country label shown on a handful of rows.
# - Updated case/death/recovered counters for several countries.
# - Updated "Datos actualizados" timestamp (18:05 -> 19:05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 19:05"
$ws.Range("B4").Value = 1803135
$ws.Range("C4").Value = 9605
$ws.Range("D4").Value = 520996
$ws.Range("E4").Value = 1177229
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = 104910
$ws.Range("B9").Value = 232664
$ws.Range("C9").Value = 416
$ws.Range("D9").Value = 155633
$ws.Range("E9").Value = 43691
$ws.Range("G9").Value = 111
$ws.Range("H9").Value = 33340
$ws.Range("B11").Value = 183149
$ws.Range("C11").Value = 130
$ws.Range("E11").Value = 9651
$ws.Range("B12").Value = 181401
$ws.Range("C12").Value = 7910
$ws.Range("D12").Value = 86668
$ws.Range("E12").Value = 89563
$ws.Range("G12").Value = 190
$ws.Range("H12").Value = 5170
$ws.Range("B17").Value = 90161
$ws.Range("C17").Value = 743
$ws.Range("D17").Value = 48050
$ws.Range("E17").Value = 35038
$ws.Range("G17").Value = 94
$ws.Range("H17").Value = 7073
$ws.Range("B37").Value = 24929
$ws.Range("C37").Value = 53
$ws.Range("E37").Value = 1189
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 1651
$ws.Range("E47").Value = 10107
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 524
$ws.Range("A92").Value = "Somalia"
$ws.Range("B92").Value = 1916
$ws.Range("C92").Value = 88
$ws.Range("D92").Value = 327
$ws.Range("E92").Value = 1516
$ws.Range("H92").Value = 73
$ws.Range("A93").Value = "Kenia"
$ws.Range("B93").Value = 1888
$ws.Range("C93").Value = 143
$ws.Range("D93").Value = 464
$ws.Range("E93").Value = 1361
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 63
$ws.Range("A94").Value = "Estonia"
$ws.Range("B94").Value = 1865
$ws.Range("C94").Value = 6
$ws.Range("D94").Value = 1622
$ws.Range("E94").Value = 176
$ws.Range("H94").Value = 67
$ws.Range("A119").Value = "Paraguay"
$ws.Range("B119").Value = 964
$ws.Range("C119").Value = 47
$ws.Range("D119").Value = 466
$ws.Range("E119").Value = 487
$ws.Range("H119").Value = 11
$ws.Range("A120").Value = "Niger"
$ws.Range("B120").Value = 955
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 813
$ws.Range("E120").Value = 78
$ws.Range("H120").Value = 64
$ws.Range("A121").Value = "Republica de Chipre"
$ws.Range("B121").Value = 944
$ws.Range("C121").Value = 2
$ws.Range("D121").Value = 784
$ws.Range("E121").Value = 143
$ws.Range("H121").Value = 17
$ws.Range("D128").Value = 470
$ws.Range("E128").Value = 224
$ws.Range("A196").Value = "Curazao"
$ws.Range("C196").Value = 1
$ws.Range("D196").Value = 14
$ws.Range("E196").Value = 4
$ws.Range("H196").Value = 1
$ws.Range("A197").Value = "Laos"
$ws.Range("D197").Value = 16
$ws.Range("E197").Value = 3
$ws.Range("A198").Value = "Nueva Caledonia"
$ws.Range("B198").Value = 19
$ws.Range("D198").Value = 18
$ws.Range("E198").Value = 1
$ws.Range("A199").Value = "Fiyi"
$ws.Range("D199").Value = 15
$ws.Range("H199").Value = 0
